$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("F4").Value = 1.59
$ws.Range("G4").Value = 1.69
$ws.Range("H4").Value = 5.9
$ws.Range("P4").Value = 1.95

# Row 5
$ws.Range("P5").Value = 1.68
$ws.Range("Q5").Value = 1.97

# Row 6
$ws.Range("F6").Value = 1.71
$ws.Range("G6").Value = 2.02
$ws.Range("H6").Value = 1.98
$ws.Range("J6").Value = 1.98
$ws.Range("P6").Value = 1.91
$ws.Range("Q6").Value = 1.72

# Row 7
$ws.Range("F7").Value = 1.57
$ws.Range("P7").Value = 1.47
$ws.Range("Q7").Value = 2.28

# Row 9
$ws.Range("G9").Value = 2.26
$ws.Range("H9").Value = 4.4
$ws.Range("I9").Value = 4.7

# Row 11
$ws.Range("F11").Value = 2.52
$ws.Range("H11").Value = 2.72
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.7
$ws.Range("K11").Value = 4.4
$ws.Range("Q11").Value = 2.12

# Row 12
$ws.Range("G12").Value = 3.25
$ws.Range("J12").Value = 2.26

# Row 13
$ws.Range("F13").Value = 2.18
$ws.Range("G13").Value = 3.05
$ws.Range("H13").Value = 1.49
$ws.Range("I13").Value = 1000
$ws.Range("J13").Value = 1.49
$ws.Range("K13").Value = 1000
$ws.Range("P13").Value = 1.53
$ws.Range("Q13").Value = 2.04

# Row 14
$ws.Range("P14").Value = 1.63
$ws.Range("Q14").Value = 1.89
